# Correction type pour génération à partir fsh ea4a6f04ed193a83290686b2f69a3f9cd2e7f4ad
#
# - "Metadata" sheet, "Name" row (row 4): fill in the ValueSet name "TypehoraireVs"
#   in column B (it was previously blank).
# - "Metadata" sheet, "Date" row (row 8): bump the generation timestamp in column B
#   to reflect the regeneration.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B4").Value = "TypehoraireVs"
$ws.Range("B8").Value = "2025-07-18T06:40:38+00:00"
